$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header test case description (B1)
$ws.Range("B1").Value = "Test Case: Testing to see any errors while logging in"

# Rewrite the existing step rows (2-4) with clarified text
$ws.Range("C2").Value = "Step 1: Sign-in with no fields filled with information"
$ws.Range("D2").Value = "An error text pops up saying that I need to fill out those field"

$ws.Range("C3").Value = "Step 2: Sign in with email filled out but no password"
$ws.Range("D3").Value = "An error text pops up saying that I need to fill out the password field"

$ws.Range("C4").Value = "Step 3: Sign in with password filled out but not email"
$ws.Range("D4").Value = "An error text pops up saying that the email is required"

# Add two new step rows (5-6) that previously were blank
$ws.Range("C5").Value = "Step 4: Fill out both fields with the wrong information"
$ws.Range("D5").Value = "An error pops up saying that the email or password is wrong"

$ws.Range("C6").Value = "Step 5: Fill out the fields with the correct information "
$ws.Range("D6").Value = "I am logged in and redirected to the user's dashboard with no errors."

# Match the same formatting used by the rest of the steps table (top aligned, wrapped text)
$ws.Range("C5:D6").WrapText = $true
$ws.Range("C5:D6").VerticalAlignment = -4160

# Update the saved selection to match the edited state
$ws.Range("D6").Select()
